# Updated cryptos list on Fri Mar 24 09:52:17 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures that were
# re-scraped from coinranking.com, and re-orders four coins (rows 37-40)
# whose relative ranking changed between scrapes:
#   - TheSandbox now ranks above InternetComputer(DFINITY)
#   - Hedera now ranks above Algorand
#
# Price cells are plain text in the source sheet (not numbers - note values
# like "28.152.21" aren't valid numeric literals), so each D-column write
# temporarily forces a Text number format before assigning the string and
# then restores the default "Normal" style so the cell's on-disk style index
# is left exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.111.94"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +1.44%  "

# --- Row 3 ---
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.791.00"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +1.78%  "

# --- Row 4 ---
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

# --- Row 5 ---
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "323.24"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.02%  "

# --- Row 6 ---
$ws.Range("E6").Value = "  +0.03%  "

# --- Row 7 ---
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4278"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -3.27%  "

# --- Row 8 ---
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3622"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -3.13%  "

# --- Row 9 ---
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "44.69"
$cell.Style = "Normal"

# --- Row 10 ---
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07539"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -2.87%  "

# --- Row 11 ---
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.113"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "

# --- Row 12 ---
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "

# --- Row 13 ---
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "21.65"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "

# --- Row 14 ---
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.145"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.67%  "

# --- Row 15 ---
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.331"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.94%  "

# --- Row 16 ---
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.813.04"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +3.14%  "

# --- Row 17 ---
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "91.83"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +2.24%  "

# --- Row 18 ---
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.00001071"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.90%  "

# --- Row 19 ---
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06343"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.95%  "

# --- Row 20 ---
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.01%  "

# --- Row 21 ---
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "17.16"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.09%  "

# --- Row 22 ---
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.975"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -3.38%  "

# --- Row 23 ---
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "28.147.20"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.45%  "

# --- Row 24 ---
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.38"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -2.30%  "

# --- Row 25 ---
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.172"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -6.16%  "

# --- Row 26 ---
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "159.05"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +3.54%  "

# --- Row 27 ---
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "20.32"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -2.20%  "

# --- Row 28 ---
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.016.67"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +3.12%  "

# --- Row 29 ---
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.226"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -5.43%  "

# --- Row 30 ---
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "127.61"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.84%  "

# --- Row 31 ---
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.168"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -3.82%  "

# --- Row 32 ---
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.825"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.99%  "

# --- Row 33 ---
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.09000"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -3.26%  "

# --- Row 34 ---
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.534"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -3.00%  "

# --- Row 35 ---
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "12.72"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.21%  "

# --- Row 36 ---
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.02354"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.33%  "

# --- Rows 37-40: coin ranking reshuffled (TheSandbox/InternetComputer swap,
#     Hedera/Algorand swap) with refreshed price/volume data ---
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.6491"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "5.085"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.06100"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.80%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.2115"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -3.34%  "


# --- Row 41 ---
$ws.Range("E41").Value = "  -0.50%  "

# --- Row 42 ---
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.427"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.63%  "

# --- Row 43 ---
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "7.928"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.86%  "

# --- Row 44 ---
$ws.Range("E44").Value = "  -0.08%  "

# --- Row 45 ---
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "13.63"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.88%  "

# --- Row 46 ---
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.6004"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.34%  "

# --- Row 47 ---
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "3.707"
$cell.Style = "Normal"

# --- Row 48 ---
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "125.13"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.94%  "

# --- Row 49 ---
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.997"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.02%  "

# --- Row 50 ---
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.150"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.56%  "

# --- Row 51 ---
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.06966"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.88%  "

